# Fixed Duplicate Cable Drawing
# Row 3 duplicated the cable-drawing info that belongs with row 2 (it only
# carried the second pull's cable size). Fold that into row 2 and drop the
# now-duplicate row, update the bundle diameter/weight to the corrected
# values, and drop the now-unused helper columns I:J along with the
# vertical 2-row merges that spanned rows 2:3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the first pull's info; the real cable size for the bundle is
# the one that was duplicated onto row 3 ("2C#2"), and the diameter/weight
# are corrected to their real (non-duplicated) values.
$ws.Range("E2").Value = "2C#2"
$ws.Range("G2").Value = 2.27
$ws.Range("H2").Value = 0.89

# Remove the now-duplicate row 3 entirely (also collapses the row 2:3
# merges back down to single cells, so Excel drops them).
$ws.Rows.Item(3).Delete()

# Columns I:J were only ever used to host the (now gone) vertical merges;
# remove them so the sheet's used range shrinks back down to H.
$ws.Range("I1:J1").EntireColumn.Delete()
